$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta...
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5b616da250>),
                ('model',
                 BaggingClassifier(estimator=LogisticRegression(class_weight='balanced',
                                                                l1_ratio=0.7,
                                                                max_iter=1000,
                                                                penalty='elasticnet',
                                                                random_state=42,
                                                                solver='saga'),
                                   n_estimators=300, random_state=42))])
'@
$ws.Range("C2").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f5b619ff8b0>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 300, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'elasticnet', 'model__estimator__l1_ratio': 0.7, 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("D2").Value = 0.624944665570591
$ws.Range("E2").Value = 'Random'
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 0.6441693125669303
$ws.Range("H2").Value = 0.4554673721340389
$ws.Range("I2").Value = '[1 1 1 0 1 0 1 0 1 0 1 1 1 0 0 0 0 1 0 0 0 0 1 0]'
$ws.Range("J2").Value = '[0 0 0 1 0 1 1 1 0 0 0 1 1 0 1 0 0 1 0 1 1 0 0 0]'

# Row 3
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta...
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5b616da1c0>),
                ('model',
                 BaggingClassifier(estimator=LogisticRegression(class_weight='balanced',
                                                                l1_ratio=0.01,
                                                                max_iter=1000,
                                                                penalty='elasticnet',
                                                                random_state=42,
                                                                solver='saga'),
                                   n_estimators=300, random_state=42))])
'@
$ws.Range("C3").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f5aea36e3d0>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 300, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'elasticnet', 'model__estimator__l1_ratio': 0.01, 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("D3").Value = 0.5700556065881143
$ws.Range("E3").Value = 'Random'
$ws.Range("F3").Value = 69
$ws.Range("G3").Value = 0.6430859674212038
$ws.Range("H3").Value = 0.529189352692075
$ws.Range("I3").Value = '[0 1 1 0 1 0 0 0 1 1 1 0 0 0 1 0 1 0 1 1 0 0 1 0]'
$ws.Range("J3").Value = '[0 1 0 0 0 0 0 0 0 0 1 1 0 0 0 1 0 0 1 0 0 1 1 1]'

# Row 4
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta...
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5b64062f10>),
                ('model',
                 BaggingClassifier(estimator=LogisticRegression(class_weight='balanced',
                                                                l1_ratio=0.7,
                                                                max_iter=1000,
                                                                penalty='elasticnet',
                                                                random_state=42,
                                                                solver='saga'),
                                   random_state=42))])
'@
$ws.Range("C4").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f5b63fefd60>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 10, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'elasticnet', 'model__estimator__l1_ratio': 0.7, 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("D4").Value = 0.6392788583965053
$ws.Range("E4").Value = 'Random'
$ws.Range("F4").Value = 23
$ws.Range("G4").Value = 0.6497870906473057
$ws.Range("H4").Value = 0.4928571428571429
$ws.Range("I4").Value = '[0 1 1 0 0 1 0 0 0 1 0 1 0 1 0 1 0 0 1 0 0 1 1 1]'
$ws.Range("J4").Value = '[0 0 0 1 0 1 0 0 0 0 1 0 0 0 1 1 0 0 0 1 1 1 0 1]'

# Row 5
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta...
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5b63fefbb0>),
                ('model',
                 BaggingClassifier(estimator=LogisticRegression(class_weight='balanced',
                                                                l1_ratio=0.95,
                                                                max_iter=1000,
                                                                penalty='elasticnet',
                                                                random_state=42,
                                                                solver='saga'),
                                   random_state=42))])
'@
$ws.Range("C5").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f5bc99a4f40>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 10, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'elasticnet', 'model__estimator__l1_ratio': 0.95, 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("D5").Value = 0.6255540064996251
$ws.Range("E5").Value = 'Random'
$ws.Range("F5").Value = 42
$ws.Range("G5").Value = 0.5076582122629894
$ws.Range("H5").Value = 0.4419413919413919
$ws.Range("I5").Value = '[1 0 1 0 0 0 1 1 1 1 1 1 1 0 0 0 0 0 0 0 1 1 0 0]'
$ws.Range("J5").Value = '[0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1]'

# Row 6
$ws.Range("A6").Value = 0
$ws.Range("B6").Value = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta...
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5bc99a4190>),
                ('model',
                 BaggingClassifier(estimator=LogisticRegression(class_weight='balanced',
                                                                l1_ratio=0.5,
                                                                max_iter=1000,
                                                                penalty='elasticnet',
                                                                random_state=42,
                                                                solver='saga'),
                                   n_estimators=5, random_state=42))])
'@
$ws.Range("C6").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f5b614faaf0>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 5, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'elasticnet', 'model__estimator__l1_ratio': 0.5, 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("D6").Value = 0.5536891712406586
$ws.Range("E6").Value = 'Random'
$ws.Range("F6").Value = 89
$ws.Range("G6").Value = 0.6112767687720249
$ws.Range("H6").Value = 0.539241622574956
$ws.Range("I6").Value = '[1 0 1 0 1 1 0 0 0 1 0 1 0 1 1 1 0 1 0 1 0 0 0 0]'
$ws.Range("J6").Value = '[0 0 1 0 1 0 0 0 1 0 0 1 0 0 0 1 0 0 0 1 1 1 1 1]'
